# level 2 stuff, level complete initial setup
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new "victory"/"VICTORY" row right after the header block
#    (before the existing "climate"/"Climate" row, currently row 9).
$ws.Rows("9:9").Insert()
$ws.Range("A9").Value = "victory"
$ws.Range("B9").Value = "VICTORY"

# 2) Insert 6 new weather rows (Rain / Overcast / Fog) before the
#    "unitAllyMallet" row. Before step 1 that row was 26; after the
#    insert above it shifted to row 27.
$ws.Rows("27:32").Insert()
$ws.Range("A27").Value = "weatherRain"
$ws.Range("B27").Value = "Rain"
$ws.Range("A28").Value = "weatherRainDesc"
$ws.Range("B28").Value = "Drips."
$ws.Range("A29").Value = "weatherOvercast"
$ws.Range("B29").Value = "Overcast"
$ws.Range("A30").Value = "weatherOvercastDesc"
$ws.Range("B30").Value = "No sun."
$ws.Range("A31").Value = "weatherFog"
$ws.Range("B31").Value = "Fog"
$ws.Range("A32").Value = "weatherFogDesc"
$ws.Range("B32").Value = "Condensed BS"

# 3) Insert a new "locationPacificNorthwest"/"Pacific Northwest" row
#    before the last "tutorialDragInstruction" row. Before the prior
#    inserts that row was 56; now (after +1 and +6) it is row 63.
$ws.Rows("63:63").Insert()
$ws.Range("A63").Value = "locationPacificNorthwest"
$ws.Range("B63").Value = "Pacific Northwest"

# Update the sheet view to match: scrolled so row 16 is at the top,
# with A27 selected.
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("A27").Select()
